$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Clear old header row content (B1:E1 originally held detail/start_date/end_date/status) ----
$ws.Range("A1:E1").ClearContents()

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 24.830729166666668
$ws.Columns.Item(5).ColumnWidth = 24.330729166666668

# ---- Row 1: big title "Hello World" ----
$ws.Range("A1").Value = "Hello World"
$ws.Range("A1:E1").Merge()
$r1 = $ws.Range("A1:E1")
$r1.Font.Size = 23
$r1.HorizontalAlignment = -4108
$r1.Interior.Pattern = 1
$r1.Interior.ThemeColor = 2
$r1.Borders.LineStyle = 1
$ws.Rows.Item(1).RowHeight = 30

# ---- Row 2: "Danh sách công việc" ----
$ws.Range("A2").Value = "Danh sách công việc"
$ws.Range("A2:E2").Merge()
$r2 = $ws.Range("A2:E2")
$r2.Font.Bold = $true
$r2.HorizontalAlignment = -4108

# ---- Row 3: "Ngày (date): " ----
$ws.Range("A3").Value = "Ngày (date): "
$r3 = $ws.Range("A3:E3")
$r3.Interior.Pattern = 1
$r3.Interior.ThemeColor = 2
$r3.VerticalAlignment = -4108

# ---- Row 4: "Nhân viên ( staff) :" ----
$ws.Range("A4").Value = "Nhân viên ( staff) :"
$r4 = $ws.Range("A4:E4")
$r4.Interior.Pattern = 1
$r4.Interior.ThemeColor = 2

# ---- Row 5: table header ----
$ws.Range("A5").Value = "id"
$ws.Range("B5").Value = "Chi tết công việc`n(detail)"
$ws.Range("C5").Value = "Ngày bắt đầu`n(start_date)"
$ws.Range("D5").Value = "Ngày kết thúc`n(end_date)"
$ws.Range("E5").Value = "Trạng thái`n(end_date)"
$r5 = $ws.Range("A5:E5")
$r5.Interior.Pattern = 1
$r5.Interior.ThemeColor = 5
$r5.HorizontalAlignment = -4108
$r5.VerticalAlignment = -4108
$r5.WrapText = $true
$r5.Borders.LineStyle = 1
$r5.Borders.Item(8).LineStyle = -4142
$ws.Rows.Item(5).RowHeight = 34

# ---- Row 6: spacer row ----
$ws.Rows.Item(6).RowHeight = 20

# ---- Row 7: G7 style only (Times New Roman font) ----
$ws.Range("G7").Font.Name = "Times New Roman"

# ---- Row 8: A8:E8 style only ----
$r8 = $ws.Range("A8:E8")
$r8.VerticalAlignment = -4108

# ---- Selection ----
$ws.Range("B4").Select()

Write-Host "done"
